$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per repull of data
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = -3
$ws.Range("F9").Value = 10
$ws.Range("F11").Value = -10
$ws.Range("F12").Value = 4
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = -5
